$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting existing rows 137:274 down to 138:275
$ws.Rows("137:137").Insert()

# Populate the newly inserted row 137 with the new record's data
$ws.Range("A137").Value2 = 10
$ws.Range("B137").Value2 = 'Vega Modelo de Temuco'
$ws.Range("C137").Value2 = 'La Araucanía'
$ws.Range("D137").Value2 = 44789
$ws.Range("E137").Value2 = 9
$ws.Range("F137").Value2 = 100112039
$ws.Range("G137").Value2 = 'Ciboulette'
$ws.Range("H137").Value2 = 'Sin especificar'
$ws.Range("I137").Value2 = 'Primera'
$ws.Range("J137").Value2 = 65
$ws.Range("K137").Value2 = 7000
$ws.Range("L137").Value2 = 7000
$ws.Range("M137").Value2 = 7000
$ws.Range("N137").Value2 = '$/docena de atados'
$ws.Range("O137").Value2 = 'Provincia de Cautín'
$ws.Range("P137").Value2 = 2333
$ws.Range("Q137").Value2 = 3
$ws.Range("R137").Value2 = 'Hortaliza'
